$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing PriceChange / UpDown values for row 9
$ws.Cells.Item(9, 24).Value = -0.21000099999999833
$ws.Cells.Item(9, 25).Value = "Down"

# Copy formats from row 9 to the new row 10 before setting values, so the
# same (existing) style indexes get reused instead of creating new ones.
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(9, 19).Copy()
$ws.Cells.Item(10, 19).PasteSpecial(-4122)

$ws.Cells.Item(9, 20).Copy()
$ws.Cells.Item(10, 20).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Add a new row of data (row 10)
$ws.Cells.Item(10, 1).Value = 42653.879479166666
$ws.Cells.Item(10, 2).Value = 13
$ws.Cells.Item(10, 3).Value = "Buy"
$ws.Cells.Item(10, 4).Value = 26
$ws.Cells.Item(10, 5).Value = 8519
$ws.Cells.Item(10, 6).Value = 865
$ws.Cells.Item(10, 7).Value = 61
$ws.Cells.Item(10, 8).Value = 36
$ws.Cells.Item(10, 9).Value = 88
$ws.Cells.Item(10, 10).Value = 11
$ws.Cells.Item(10, 11).Value = 5940
$ws.Cells.Item(10, 12).Value = 151
$ws.Cells.Item(10, 13).Value = 88
$ws.Cells.Item(10, 14).Value = 47
$ws.Cells.Item(10, 15).Value = 6
$ws.Cells.Item(10, 16).Value = "Bag"
$ws.Cells.Item(10, 17).Value = 47.96375473473072
$ws.Cells.Item(10, 18).Value = 0.49
$ws.Cells.Item(10, 19).Value = 0.0521
$ws.Cells.Item(10, 20).Value = -0.0214
$ws.Cells.Item(10, 21).Value = 2.25
$ws.Cells.Item(10, 22).Value = "N/A"
$ws.Cells.Item(10, 23).Value = 0

$wb.Save()
